$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.710608243942261
$ws.Range("B1").Value = 1.732973694801331
$ws.Range("C1").Value = 1.577597260475159
$ws.Range("D1").Value = 1.201207876205444
$ws.Range("E1").Value = 0.7166420817375183
